$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-45 (timestamp / IGCC Import (MW) / IGCC Export (MW))
# replacing the previous 2025-... data (rows 2-32) with retrained/new data
# (rows 2-45), per "Retraining the model for Horeco".
$rows = @(
    @(2, 46048, 9.788, 0.275),
    @(3, 46048.01041666666, 37.428, 0),
    @(4, 46048.02083333334, 37.133, 0),
    @(5, 46048.03125, 16.765, 0),
    @(6, 46048.04166666666, 2.604, 0.004),
    @(7, 46048.05208333334, 5.466, 0.007),
    @(8, 46048.0625, 6.062, 0.068),
    @(9, 46048.07291666666, 20.84, 0),
    @(10, 46048.08333333334, 4.228, 0.965),
    @(11, 46048.09375, 0.681, 1.638),
    @(12, 46048.10416666666, 3.979, 0.546),
    @(13, 46048.11458333334, 0, 0.503),
    @(14, 46048.125, 0, 0.594),
    @(15, 46048.13541666666, 0, 1.625),
    @(16, 46048.14583333334, 0.025, 0.013),
    @(17, 46048.15625, 4.266, 0.032),
    @(18, 46048.16666666666, 21.073, 0),
    @(19, 46048.17708333334, 3.969, 0.02),
    @(20, 46048.1875, 15.651, 0),
    @(21, 46048.19791666666, 20.135, 0),
    @(22, 46048.20833333334, 28.212, 0),
    @(23, 46048.21875, 13.501, 0.108),
    @(24, 46048.22916666666, 0, 2.47),
    @(25, 46048.23958333334, 5.02, 0.477),
    @(26, 46048.25, 0, 14.661),
    @(27, 46048.26041666666, 0, 20.623),
    @(28, 46048.27083333334, 0, 70.763),
    @(29, 46048.28125, 0, 50.94),
    @(30, 46048.29166666666, 0, 43.798),
    @(31, 46048.30208333334, 0, 29.032),
    @(32, 46048.3125, 0.034, 3.997),
    @(33, 46048.32291666666, 4.742, 0.81),
    @(34, 46048.33333333334, 0.706, 4.159),
    @(35, 46048.34375, 2.082, 0),
    @(36, 46048.35416666666, 6.09, 0),
    @(37, 46048.36458333334, 6.737, 0),
    @(38, 46048.375, 36.113, 0),
    @(39, 46048.38541666666, 26.608, 0),
    @(40, 46048.39583333334, 15.146, 0),
    @(41, 46048.40625, 3.547, 2.838),
    @(42, 46048.41666666666, 0, 7.613),
    @(43, 46048.42708333334, 0, 19.572),
    @(44, 46048.4375, 8.072, 0.286),
    @(45, 46048.44791666666, 0, 0)
)

# Remember the number format used on the existing timestamp column (A2:A32)
# so it can be (re)applied to the newly created rows (33:45) in column A.
$dateFormat = $ws.Range("A2").NumberFormat

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
}

# Ensure newly added rows (33:45) carry the same timestamp number format as
# the rest of column A.
$ws.Range("A33:A45").NumberFormat = $dateFormat
